$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 333337340
$ws.Range("I94").Value = 2000
$ws.Range("J94").Value = 500005000
$ws.Range("K94").Value = 2000
$ws.Range("L94").Value = 500005000
$ws.Range("M94").Value = -1549
$ws.Range("N94").Value = -500005902
$ws.Range("H125").Value = 1807.1666
$ws.Range("I125").Value = 2382.6667
$ws.Range("J125").Value = 1519.4166
$ws.Range("K125").Value = 21444.0003
$ws.Range("L125").Value = 13674.7494
$ws.Range("M125").Value = -18984.0003
$ws.Range("N125").Value = -18594.7494
$ws.Range("H129").Value = 800.78125
$ws.Range("I129").Value = 512.1429000000001
$ws.Range("J129").Value = 881.6
$ws.Range("K129").Value = 1536.4287
$ws.Range("L129").Value = 2644.8
$ws.Range("M129").Value = 3463.5713
$ws.Range("N129").Value = -12644.8
$ws.Range("H132").Value = 5685559
$ws.Range("I132").Value = 5817758
$ws.Range("K132").Value = 17453274
$ws.Range("M132").Value = -17450744
$ws.Range("H135").Value = 564.7778
$ws.Range("I135").Value = 617.7619
$ws.Range("K135").Value = 5559.857099999999
$ws.Range("M135").Value = -3024.857099999999
$ws.Range("H137").Value = 1671.2354
$ws.Range("I137").Value = 1084.25
$ws.Range("J137").Value = 3080
$ws.Range("K137").Value = 3252.75
$ws.Range("L137").Value = 9240
$ws.Range("M137").Value = -702.75
$ws.Range("N137").Value = -14340
$ws.Range("H138").Value = 4482.24
$ws.Range("I138").Value = 1662.9048
$ws.Range("J138").Value = 6523.8276
$ws.Range("K138").Value = 4988.7144
$ws.Range("L138").Value = 19571.4828
$ws.Range("M138").Value = 151.2856000000002
$ws.Range("N138").Value = -29851.4828

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2988.2778
$ws.Range("I61").Value = 1196
$ws.Range("J61").Value = 3212.3125
$ws.Range("K61").Value = 1196
$ws.Range("L61").Value = 3212.3125
$ws.Range("M61").Value = -984
$ws.Range("N61").Value = -3636.3125
$ws.Range("H132").Value = 2998.121
$ws.Range("I132").Value = 2914.8147
$ws.Range("J132").Value = 3373
$ws.Range("K132").Value = 8744.444100000001
$ws.Range("L132").Value = 10119
$ws.Range("M132").Value = -6214.444100000001
$ws.Range("N132").Value = -15179
$ws.Range("H134").Value = 64836.668
$ws.Range("J134").Value = 64836.668
$ws.Range("L134").Value = 64836.668
$ws.Range("N134").Value = -74976.66800000001
$ws.Range("H135").Value = 57450
$ws.Range("J135").Value = 57450
$ws.Range("L135").Value = 57450
$ws.Range("N135").Value = -67590
$ws.Range("H136").Value = 2988.2778
$ws.Range("I136").Value = 1196
$ws.Range("J136").Value = 3212.3125
$ws.Range("K136").Value = 3588
$ws.Range("L136").Value = 9636.9375
$ws.Range("M136").Value = -1038
$ws.Range("N136").Value = -14736.9375
$ws.Range("H137").Value = 65000
$ws.Range("J137").Value = 65000
$ws.Range("L137").Value = 65000
$ws.Range("N137").Value = -75200
$ws.Range("H139").Value = 51678.75
$ws.Range("J139").Value = 51678.75
$ws.Range("L139").Value = 51678.75
$ws.Range("N139").Value = -61958.75
$ws.Range("H141").Value = 30390
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2009
$ws.Range("I134").Value = 2104.087
$ws.Range("J134").Value = 1462.25
$ws.Range("K134").Value = 6312.261
$ws.Range("L134").Value = 4386.75
$ws.Range("M134").Value = -3777.261
$ws.Range("N134").Value = -9456.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 6425.3516
$ws.Range("I5").Value = 891.8461
$ws.Range("J5").Value = 19504.545
$ws.Range("K5").Value = 2675.5383
$ws.Range("L5").Value = 58513.63499999999
$ws.Range("M5").Value = -2563.5383
$ws.Range("N5").Value = -58737.63499999999
$ws.Range("H131").Value = 854.66
$ws.Range("J131").Value = 877.87915
$ws.Range("L131").Value = 2633.63745
$ws.Range("N131").Value = -12713.63745
$ws.Range("H133").Value = 6442.857
$ws.Range("J133").Value = 7350
$ws.Range("L133").Value = 22050
$ws.Range("N133").Value = -32170
$ws.Range("H134").Value = 4004.2
$ws.Range("I134").Value = 4004.2
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 12012.6
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -6942.599999999999
$ws.Range("N134").Value = $null
$ws.Range("H135").Value = 6425.3516
$ws.Range("I135").Value = 891.8461
$ws.Range("J135").Value = 19504.545
$ws.Range("K135").Value = 8026.6149
$ws.Range("L135").Value = 175540.905
$ws.Range("M135").Value = -5491.6149
$ws.Range("N135").Value = -180610.905
$ws.Range("H136").Value = 1607.2
$ws.Range("I136").Value = 1109
$ws.Range("J136").Value = 3600
$ws.Range("K136").Value = 3327
$ws.Range("L136").Value = 10800
$ws.Range("M136").Value = 1773
$ws.Range("N136").Value = -21000
$ws.Range("H137").Value = 16673947
$ws.Range("I137").Value = 8736
$ws.Range("J137").Value = 100000000
$ws.Range("K137").Value = 26208
$ws.Range("L137").Value = 300000000
$ws.Range("M137").Value = -21108
$ws.Range("N137").Value = -300010200
$ws.Range("H139").Value = 2381
$ws.Range("I139").Value = 1512.3846
$ws.Range("J139").Value = 3994.1428
$ws.Range("K139").Value = 4537.1538
$ws.Range("L139").Value = 11982.4284
$ws.Range("M139").Value = 602.8462
$ws.Range("N139").Value = -22262.4284
$ws.Range("H140").Value = 1771.1875
$ws.Range("I140").Value = 1308.9
$ws.Range("J140").Value = 2541.6667
$ws.Range("K140").Value = 3926.7
$ws.Range("L140").Value = 7625.000100000001
$ws.Range("M140").Value = 1253.3
$ws.Range("N140").Value = -17985.0001
$ws.Range("H141").Value = 3983.3333
$ws.Range("I141").Value = 3500
$ws.Range("J141").Value = 4950
$ws.Range("K141").Value = 10500
$ws.Range("L141").Value = 14850
$ws.Range("M141").Value = -5320
$ws.Range("N141").Value = -25210

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 232775.53
$ws.Range("I102").Value = 2102.7827
$ws.Range("J102").Value = 2001266.6
$ws.Range("K102").Value = 2102.7827
$ws.Range("L102").Value = 2001266.6
$ws.Range("M102").Value = -480.7827000000002
$ws.Range("N102").Value = -2004510.6
$ws.Range("H132").Value = 2496.923
$ws.Range("I132").Value = 2151.8333
$ws.Range("J132").Value = 2792.7144
$ws.Range("K132").Value = 6455.499899999999
$ws.Range("L132").Value = 8378.143199999999
$ws.Range("M132").Value = -3925.499899999999
$ws.Range("N132").Value = -13438.1432
$ws.Range("H135").Value = 27204.727
$ws.Range("J135").Value = 27204.727
$ws.Range("L135").Value = 27204.727
$ws.Range("N135").Value = -37344.727
$ws.Range("H137").Value = 68000
$ws.Range("J137").Value = 68000
$ws.Range("L137").Value = 68000
$ws.Range("N137").Value = -78200
$ws.Range("H138").Value = 59962.25
$ws.Range("I138").Value = 54500
$ws.Range("J138").Value = 65424.5
$ws.Range("K138").Value = 54500
$ws.Range("L138").Value = 65424.5
$ws.Range("M138").Value = -49360
$ws.Range("N138").Value = -75704.5
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").Value = $null
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1832.8235
$ws.Range("I136").Value = 1767
$ws.Range("J136").Value = 1906.875
$ws.Range("K136").Value = 5301
$ws.Range("L136").Value = 5720.625
$ws.Range("M136").Value = -2751
$ws.Range("N136").Value = -10820.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 10000
$ws.Range("J41").Value = 10000
$ws.Range("L41").Value = 10000
$ws.Range("N41").Value = -10780
$ws.Range("H45").Value = 8951.666999999999
$ws.Range("J45").Value = 8951.666999999999
$ws.Range("L45").Value = 8951.666999999999
$ws.Range("N45").Value = -9933.666999999999
$ws.Range("H74").Value = 9813
$ws.Range("J74").Value = 9813
$ws.Range("L74").Value = 9813
$ws.Range("N74").Value = -11685
$ws.Range("H77").Value = 9813
$ws.Range("J77").Value = 9813
$ws.Range("L77").Value = 29439
$ws.Range("N77").Value = -38799
$ws.Range("H81").Value = 333963
$ws.Range("I81").Value = 500470
$ws.Range("J81").Value = 250709.5
$ws.Range("K81").Value = 1000940
$ws.Range("L81").Value = 501419
$ws.Range("M81").Value = -999879
$ws.Range("N81").Value = -503541
$ws.Range("H84").Value = 333963
$ws.Range("I84").Value = 500470
$ws.Range("J84").Value = 250709.5
$ws.Range("K84").Value = 5004700
$ws.Range("L84").Value = 2507095
$ws.Range("M84").Value = -4999396
$ws.Range("N84").Value = -2517703
$ws.Range("H113").Value = 818.9167
$ws.Range("I113").Value = 542.25
$ws.Range("J113").Value = 1372.25
$ws.Range("K113").Value = 1626.75
$ws.Range("L113").Value = 4116.75
$ws.Range("M113").Value = 543.25
$ws.Range("N113").Value = -8456.75
$ws.Range("H126").Value = 1208.7273
$ws.Range("I126").Value = 1168
$ws.Range("K126").Value = 3504
$ws.Range("M126").Value = -1034

Write-Output "Edits applied"